$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @{ B=3.182878228561681; C=1.65323645889881;  D=0.7127328510149897; E=0.4998867070740569; G=6.048734245549538 }
  3  = @{ B=3.182878228561681; C=1.65323645889881;  D=0.1529057820181812; E=6.48142807727062;    G=11.47044854674929 }
  4  = @{ B=3.182878228561681; C=1.65323645889881;  D=0.7127328510149897; E=6.48142807727062;    G=12.0302756157461 }
  5  = @{ B=0.3464964993005633; C=0.3375848360084654; D=0.1529057820181812; E=0.4998867070740569; G=1.336873824401267 }
  6  = @{ B=3.182878228561681; C=1.65323645889881;  D=3.082599426703578; E=6.48142807727062;    G=14.40014219143469 }
  7  = @{ B=3.182878228561681; C=1.65323645889881;  D=0.1529057820181812; E=0.4998867070740569; G=5.488907176552729 }
  8  = @{ B=0.7287194209349384; C=1.65323645889881;  D=0.7127328510149897; E=6.48142807727062;    G=9.576116808119359 }
  9  = @{ B=3.182878228561681; C=1.65323645889881;  D=3.082599426703578; E=6.48142807727062;    G=14.40014219143469 }
  10 = @{ B=0.7287194209349384; C=1.65323645889881;  D=0.7127328510149897; E=0.4998867070740569; G=3.594575437922795 }
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  $ws.Range("B$row").Value = $vals.B
  $ws.Range("C$row").Value = $vals.C
  $ws.Range("D$row").Value = $vals.D
  $ws.Range("E$row").Value = $vals.E
  $ws.Range("G$row").Value = $vals.G
}
